$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# B11 held the text "R40"; the edit replaces its content with the text "1".
# A leading apostrophe forces Excel to store the numeral as text (so it
# keeps landing in the shared-string table, t="s") instead of silently
# re-interpreting it as the number 1.
$ws.Range("B11").Value = "'1"
